$d = $word.ActiveDocument

# --- Step 1: locate "Collins" and retype it as "Berlin Heights, Ohio " ---
$r1 = $d.Content
$found = $r1.Find.Execute("Collins", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Collins' in document"
}
$r1.Text = "Berlin Heights, Ohio "

# --- Step 2: locate the "_GoBack" bookmark that now sits right after that text ---
$bm = $d.Bookmarks.Item("_GoBack")

# Find the paragraph that contains the bookmark (avoids the Range.Paragraphs quirk).
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($bm.Start -ge $p.Range.Start -and $bm.Start -le $p.Range.End) {
        $para = $p
        break
    }
}
$paraEnd = $para.Range.End - 1   # exclude the paragraph mark

# --- Step 3: the remainder of the old text (", Ohio 44814") follows the bookmark; retype as "44814" ---
$r2 = $d.Range($bm.End, $paraEnd)
$r2.Text = "44814"

# --- Step 4: re-stake the _GoBack bookmark so it wraps the freshly (re)typed "44814" ---
$bm2 = $d.Bookmarks.Item("_GoBack")
$newEnd = $para.Range.End - 1
$rngNew = $d.Range($bm2.Start, $newEnd)
$d.Bookmarks.Add("_GoBack", $rngNew)
